# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation"
#   on the Overview, zh-cn and de-de sheets.
# - Narrow the now-shorter "Status" columns (E/F on Overview, C on the
#   language sheets) to match the new, shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

# --- Update the status values -------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("C3").Value = "In Translation"

$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("C3").Value = "In Translation"

# --- Narrow the affected columns to fit the shorter text -----------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F

$wsZh.Columns.Item(3).ColumnWidth = 12.5         # column C
$wsDe.Columns.Item(3).ColumnWidth = 12.5         # column C
